# "Add padding to icons"
#
# The icon is built from a group shape ("Group 3") whose hit-test/bounding
# rectangle ("Rectangle 1", the last child) is enlarged by 10% (457200 EMU ->
# 502920 EMU) and re-centered, which in turn grows the group's own visual
# bounding box by the same amount while keeping the visible glyph (the cross
# + connector lines that make up the actual icon artwork) exactly where it
# was.
#
# PowerPoint stores shape geometry in EMU (914400 EMU/inch, 12700 EMU/point)
# but the COM object model only exposes Left/Top/Width/Height in points as
# (single-precision) floats. To land on an exact target EMU value we search
# for a point value whose float32 representation, multiplied by 12700 and
# truncated, reproduces that exact integer EMU.

function ToEmuExact($targetEmu) {
    $base = $targetEmu / 12700.0
    for ($i = 0; $i -lt 20000; $i++) {
        $delta = $i * 0.0000001
        foreach ($cand in @($base + $delta, $base - $delta)) {
            $f32 = [float]$cand
            $backEmu = [int64]([double]$f32 * 12700.0)
            if ($backEmu -eq $targetEmu) {
                return $cand
            }
        }
    }
    # Fall back to the naive conversion if an exact hit wasn't found
    # (not expected for the values used below).
    return $base
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "Group 3" is the icon's top-level group shape on the slide.
$grp = $s.Shapes.Item(9)
if ($grp.Name -ne "Group 3") {
    foreach ($sh in $s.Shapes) {
        if ($sh.Name -eq "Group 3") { $grp = $sh }
    }
}

# "Rectangle 1" is the last item inside the group - the invisible
# click-target rectangle that is being padded.
$rect = $grp.GroupItems.Item($grp.GroupItems.Count)
if ($rect.Name -ne "Rectangle 1") {
    for ($i = 1; $i -le $grp.GroupItems.Count; $i++) {
        if ($grp.GroupItems.Item($i).Name -eq "Rectangle 1") {
            $rect = $grp.GroupItems.Item($i)
        }
    }
}

# New group bounding box (slide coordinates, EMU).
$groupOffX = ToEmuExact 2308058
$groupOffY = ToEmuExact 1140167
$groupExtCx = ToEmuExact 502920
$groupExtCy = ToEmuExact 502920

# New "Rectangle 1" geometry (slide coordinates, EMU) - grown by 10% and
# centered on its former center.
$rectOffX = ToEmuExact 2345677
$rectOffY = ToEmuExact 1140167
$rectExtCx = ToEmuExact 502920
$rectExtCy = ToEmuExact 502920

$grp.Left = $groupOffX
$grp.Top = $groupOffY
$grp.Width = $groupExtCx
$grp.Height = $groupExtCy

$rect.Left = $rectOffX
$rect.Top = $rectOffY
$rect.Width = $rectExtCx
$rect.Height = $rectExtCy

Write-Output ("Group 3: Left=" + $grp.Left + " Top=" + $grp.Top + " Width=" + $grp.Width + " Height=" + $grp.Height)
Write-Output ("Rectangle 1: Left=" + $rect.Left + " Top=" + $rect.Top + " Width=" + $rect.Width + " Height=" + $rect.Height)
